$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 74
$ws.Range("I2").Value = 224
$ws.Range("J2").Value = 869
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 218
$ws.Range("M2").Value = 14
$ws.Range("N2").Value = 162
$ws.Range("P2").Value = 3
$ws.Range("R2").Value = 9
$ws.Range("S2").Value = 114
$ws.Range("T2").Value = 159
$ws.Range("U2").Value = 6
$ws.Range("V2").Value = 1327
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 1314
$ws.Range("Z2").Value = 6
$ws.Range("AA2").Value = 9
